# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the first data row
# on both the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 14:53:22"
$wsZhCn.Range("H2").Value = "2016-03-19 14:53:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 14:53:26"
$wsDeDe.Range("H2").Value = "2016-03-19 14:53:44"
